$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# New Key/Value rows to append for the options modal + hud art strings
$rows = @(
    @("options", "OPTIONS"),
    @("sound", "SOUND"),
    @("speech", "SPEECH"),
    @("on", "ON"),
    @("off", "OFF"),
    @("close", "CLOSE"),
    @("chain_of_custody", "Chain of Custody"),
    @("activity_log", "Activity Log"),
    @("help", "Help"),
    @("disk_clone", "Disk Clone"),
    @("hard_disk_drive", "Hard Disk Drive"),
    @("usb_flash_drive", "USB Flash Drive")
)

$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

$lastRow = $startRow + $rows.Count - 1
$nextRow = $lastRow + 1
$ws.Range("A$nextRow").Select()
